$wb = $excel.ActiveWorkbook
$piSheet = $wb.Worksheets("PI hours")

# --- Add new "cfop" column (G) to "PI hours" sheet ---
$piSheet.Range("F1").Copy()
$piSheet.Range("G1").PasteSpecial(-4122)
$piSheet.Range("G1").Value = "cfop"

$piSheet.Range("F2:F4").Copy()
$piSheet.Range("G2:G4").PasteSpecial(-4122)
$piSheet.Range("G2").Value = "['cfop_NH']"
$piSheet.Range("G3").Value = "['cfop_GC']"
$piSheet.Range("G4").Value = "['cfop_HUTCHINSON']"

# --- Add new "cfop hours" worksheet after "unit(accumulative) hours" ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "cfop hours"

# Header row (reuse header formatting from department hours sheet)
$deptSheet = $wb.Worksheets("department hours")
$deptSheet.Range("B1:D1").Copy()
$newSheet.Range("B1:D1").PasteSpecial(-4122)
$newSheet.Range("B1").Value = "cfop"
$newSheet.Range("C1").Value = "hours"
$newSheet.Range("D1").Value = "percentage"

# Data rows (reuse formatting from department hours sheet body)
$deptSheet.Range("A2:D4").Copy()
$newSheet.Range("A2:D4").PasteSpecial(-4122)

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "cfop_NH"
$newSheet.Range("C2").Value = 39
$newSheet.Range("D2").Value = 54.92957746478874

$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "cfop_GC"
$newSheet.Range("C3").Value = 16
$newSheet.Range("D3").Value = 22.53521126760563

$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "cfop_HUTCHINSON"
$newSheet.Range("C4").Value = 16
$newSheet.Range("D4").Value = 22.53521126760563

# Restore the original active sheet/tab selection
$piSheet.Activate()
